# Weekly update: insert a new latest-week record at the top of the data
# (row 4, just below the header + two already-current rows) and push the
# existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:18 down to 5:19, carrying their values/styles with them.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range('A4').Value2 = 4
$ws.Range('B4').Value2 = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C4').Value2 = 'Los Lagos'
$ws.Range('D4').Value2 = 44707
$ws.Range('E4').Value2 = 10
$ws.Range('F4').Value2 = 100112012
$ws.Range('G4').Value2 = 'Espinaca'
$ws.Range('H4').Value2 = 'Sin especificar'
$ws.Range('I4').Value2 = 'Primera'
$ws.Range('J4').Value2 = 15
$ws.Range('K4').Value2 = 12000
$ws.Range('L4').Value2 = 12000
$ws.Range('M4').Value2 = 12000
$ws.Range('N4').Value2 = '$/cuna 10 kilos'
$ws.Range('O4').Value2 = 'Región Metropolitana'
$ws.Range('P4').Value2 = 1200
$ws.Range('Q4').Value2 = 10
$ws.Range('R4').Value2 = 'Hortaliza'
